## Add a new worksheet "test" after "demo" and populate it with data that
## mirrors the "demo" sheet's layout (same header pattern, same per-cell
## styling), while leaving "demo" itself mostly intact apart from a
## refreshed selection/column width.

$wb = $excel.ActiveWorkbook
$demo = $wb.Worksheets.Item("demo")

# --- demo sheet: widen column A and move the selection -------------------
$demo.Columns.Item(1).ColumnWidth = 18.86
$demo.Range("H12").Select() | Out-Null

# --- create the new "test" sheet right after "demo" -----------------------
$new = $wb.Worksheets.Add($null, $demo)
$new.Name = "test"

# Reuse demo's existing cell styles (style index 1 for the A:F block, the
# default style for the G "!"/formula column) by copying formatting only -
# this keeps styles.xml untouched, exactly like the source data.
$demo.Range("A1:F2").Copy() | Out-Null
$new.Range("A1:F9").PasteSpecial(-4122) | Out-Null

$demo.Range("G2").Copy() | Out-Null
$new.Range("G2:G9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- header rows ------------------------------------------------------
$new.Range("A1").Value = "!"

$new.Range("A2").Value = "#"
$new.Range("B2").Value = "id"
$new.Range("C2").Value = "name"
$new.Range("D2").Value = "value"
$new.Range("E2").Value = "price"
$new.Range("F2").Value = "max"
$new.Range("G2").Value = "!"

# --- data rows ----------------------------------------------------------
$new.Range("B3").Value = 234234
$new.Range("C3").Value = "Della"
$new.Range("D3").Value = 100
$new.Range("E3").Value = 123
$new.Range("F3").Value = 100
$new.Range("G3").Formula = "=D3/100"

$new.Range("B4").Value = 4582
$new.Range("C4").Value = "Neka"
$new.Range("D4").Value = 200
$new.Range("E4").Value = 321
$new.Range("F4").Value = 50
$new.Range("G4").Formula = "=D4/100"

$new.Range("B5").Value = 45672
$new.Range("C5").Value = "Ares"
$new.Range("D5").Value = 300
$new.Range("E5").Value = 456
$new.Range("F5").Value = 99
$new.Range("G5").Formula = "=D5/100"

$new.Range("B6").Value = 5428
$new.Range("C6").Value = "Lome"
$new.Range("D6").Value = 400
$new.Range("E6").Value = 4895
$new.Range("F6").Value = 200
$new.Range("G6").Formula = "=D6/100"

$new.Range("B7").Value = 85
$new.Range("C7").Value = "Chopper"
$new.Range("D7").Value = 500
$new.Range("E7").Value = 21546
$new.Range("F7").Value = 200
$new.Range("G7").Formula = "=D7/100"

$new.Range("B8").Value = 62654
$new.Range("C8").Value = "Spy"
$new.Range("D8").Value = 600
$new.Range("E8").Value = 45.54
$new.Range("F8").Value = 99
$new.Range("G8").Formula = "=D8/100"

$new.Range("B9").Value = 943452
$new.Range("C9").Value = "Wolly"
$new.Range("D9").Value = 700
$new.Range("E9").Value = 453.54
$new.Range("F9").Value = 50
$new.Range("G9").Formula = "=D9/100"

# --- page setup (mirrors the rest of the workbook's print settings) -----
$ps = $new.PageSetup
$ps.PaperSize = 9
$ps.Zoom = 100
$ps.Orientation = 1
$ps.CenterHeader = '&"Times New Roman,標準"&12&A'
$ps.CenterFooter = '&"Times New Roman,標準"&12頁 &P'
$ps.LeftMargin = 56.7
$ps.RightMargin = 56.7
$ps.TopMargin = 75.8
$ps.BottomMargin = 75.8
$ps.HeaderMargin = 56.7
$ps.FooterMargin = 56.7

# --- finally, select D15 and make "test" the active sheet/tab -----------
$new.Range("D15").Select() | Out-Null
